$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 07:52"

# 2. Update Lituania figures (row 71)
$ws.Cells.Item(71, 2).Value = 1062
$ws.Cells.Item(71, 3).Value = 9
$ws.Cells.Item(71, 4).Value = 101
$ws.Cells.Item(71, 5).Value = 937
$ws.Cells.Item(71, 6).Value = 14
$ws.Cells.Item(71, 7).Value = 1
$ws.Cells.Item(71, 8).Value = 24

# 3. Update Bulgaria figures (row 82)
$ws.Cells.Item(82, 2).Value = 676
$ws.Cells.Item(82, 3).Value = 1
$ws.Cells.Item(82, 4).Value = 71
$ws.Cells.Item(82, 5).Value = 574
$ws.Cells.Item(82, 6).Value = 36
$ws.Cells.Item(82, 7).Value = 2
$ws.Cells.Item(82, 8).Value = 31

# 4. El Salvador's case count rises above Brunei's, so it moves up one row,
#    pushing Brunei / Gibraltar / Ruanda down a row each (rows 126-129).
$ws.Cells.Item(126, 1).Value = "El Salvador"
$ws.Cells.Item(126, 2).Value = 137
$ws.Cells.Item(126, 3).Value = 12
$ws.Cells.Item(126, 4).Value = 22
$ws.Cells.Item(126, 5).Value = 109
$ws.Cells.Item(126, 6).Value = 3
$ws.Cells.Item(126, 7).Value = 0
$ws.Cells.Item(126, 8).Value = 6

$ws.Cells.Item(127, 1).Value = "Brunei"
$ws.Cells.Item(127, 2).Value = 136
$ws.Cells.Item(127, 3).Value = 0
$ws.Cells.Item(127, 4).Value = 106
$ws.Cells.Item(127, 5).Value = 29
$ws.Cells.Item(127, 6).Value = 2
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 1

$ws.Cells.Item(128, 1).Value = "Gibraltar"
$ws.Cells.Item(128, 2).Value = 129
$ws.Cells.Item(128, 3).Value = 0
$ws.Cells.Item(128, 4).Value = 84
$ws.Cells.Item(128, 5).Value = 45
$ws.Cells.Item(128, 6).Value = 1
$ws.Cells.Item(128, 7).Value = 0
$ws.Cells.Item(128, 8).Value = 0

$ws.Cells.Item(129, 1).Value = "Ruanda"
$ws.Cells.Item(129, 2).Value = 126
$ws.Cells.Item(129, 3).Value = 0
$ws.Cells.Item(129, 4).Value = 25
$ws.Cells.Item(129, 5).Value = 101
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 0
